# Sequence diagram edit: method rename
#   createDefaultFileIfInvalidDate()  ->  createDefaultFileIfInvalidDateOrRange()
# The wider label text pushed the label boxes in this "lane" of the
# sequence diagram slightly left/up, so the four sibling method-label
# textboxes get nudged and the renamed one also grows wider.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# TextBox 59 - "createDefaultFileIfNotExist()" - position nudge only
$shape60 = Get-ShapeById $s 60
$shape60.Left = 29.053700787401574
$shape60.Top = 138.11275590551182

# TextBox 29 - "createDefaultFileIfSizeDiff()" - position nudge only
$shape30 = Get-ShapeById $s 30
$shape30.Left = 24.37771653543307
$shape30.Top = 172.36740157480315

# TextBox 41 - "createDefaultFileIfNull()" - position nudge only
$shape42 = Get-ShapeById $s 42
$shape42.Left = 7.903622047244094
$shape42.Top = 205.04700787401575

# TextBox 46 - renamed method, also moves and widens to fit the new text
$shape47 = Get-ShapeById $s 47
$shape47.Left = 26.276535433070865
$shape47.Top = 235.4048031496063
$shape47.Width = 198.06748031496062

$tr = $shape47.TextFrame.TextRange
$tr.Characters(1, 30).Text = "createDefaultFileIfInvalidDateOrRange"
